# NFL prediction tool: update a few input scores and select D5.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Direct input cells (C = Fav line, D = ov/un) that changed for KC/ATL/GB rows.
$ws.Range("C3").Value = -1   # KC  Fav: -3 -> -1
$ws.Range("D4").Value = 37   # ATL ov/un: 38 -> 37
$ws.Range("D5").Value = 39   # GB  ov/un: 40 -> 39

# J5 is a cached "Out:" data-table result cell (What-If data table J2:J5);
# update its stored value directly to match the new table output.
$ws.Range("J5").Value = 17   # GB  Out: 20 -> 17

# Leave the selection on D5, matching the saved sheet view.
$ws.Range("D5").Select()
